$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D (Price) and E (Volume 1h) to be treated as text so that
# numeric-looking values (e.g. "278.00", "0.95%") keep their exact original
# formatting instead of being normalized into numbers/percentages by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '278.00'
$ws.Range("E2").Value = '0.95%'
$ws.Range("D3").Value = '27.24'
$ws.Range("E3").Value = '2.60%'
$ws.Range("D4").Value = '4.873'
$ws.Range("E4").Value = '-0.36%'
$ws.Range("D5").Value = '0.06430'
$ws.Range("E5").Value = '1.31%'
$ws.Range("D6").Value = '6.992'
$ws.Range("E6").Value = '1.21%'
$ws.Range("E7").Value = '-5.95%'
$ws.Range("D8").Value = '0.8849'
$ws.Range("E8").Value = '1.96%'
$ws.Range("D9").Value = '0.1559'
$ws.Range("E9").Value = '1.40%'
$ws.Range("D10").Value = '0.05112'
$ws.Range("E10").Value = '1.98%'
$ws.Range("D11").Value = '0.07479'
$ws.Range("E11").Value = '0.92%'
$ws.Range("D12").Value = '0.02887'
$ws.Range("E12").Value = '-3.03%'
$ws.Range("D13").Value = '0.08978'
$ws.Range("E13").Value = '-0.75%'
$ws.Range("D14").Value = '0.001570'
$ws.Range("E14").Value = '-0.27%'
$ws.Range("D15").Value = '0.0006399'
$ws.Range("E15").Value = '1.27%'
$ws.Range("D16").Value = '0.006154'
$ws.Range("E16").Value = '3.95%'
$ws.Range("E17").Value = '1.07%'
$ws.Range("D18").Value = '3.314'
$ws.Range("E18").Value = '0.06%'
$ws.Range("D21").Value = '0.1350'
$ws.Range("D22").Value = '3.901'
$ws.Range("E22").Value = '0.13%'
$ws.Range("D23").Value = '0.04409'
$ws.Range("E23").Value = '1.10%'
$ws.Range("D24").Value = '0.1500'
$ws.Range("E24").Value = '8.66%'
$ws.Range("B25").Value = 'BitKan'
$ws.Range("C25").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D25").Value = '0.001175'
$ws.Range("E25").Value = '-0.05%'
$ws.Range("B26").Value = 'HotbitToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D26").Value = '0.003878'
$ws.Range("E26").Value = '-7.87%'
$ws.Range("B27").Value = 'AAXToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/LNePqkIhk+aaxtoken-aab'
$ws.Range("D27").Value = '0.08000'
$ws.Range("E27").Value = '-60.00%'
$ws.Range("E28").Value = '-1.63%'
$ws.Range("D29").Value = '0.0001650'
$ws.Range("E29").Value = '-1.76%'
$ws.Range("D40").Value = '0.04152'
$ws.Range("E40").Value = '1.15%'
$ws.Range("D41").Value = '0.006816'
$ws.Range("E41").Value = '-2.69%'
$ws.Range("E42").Value = '0.43%'
$ws.Range("D43").Value = '0.001919'
$ws.Range("E43").Value = '-10.47%'
$ws.Range("D44").Value = '0.01176'
$ws.Range("E44").Value = '9.07%'
$ws.Range("D45").Value = '0.00005305'
$ws.Range("E45").Value = '0.51%'
$ws.Range("D46").Value = '1.685'
$ws.Range("E46").Value = '13.33%'
$ws.Range("D47").Value = '0.01852'
$ws.Range("E47").Value = '-7.43%'
